# Auto-generated edit script: updates cryptos list values per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.659.83'
$ws.Range("E2").Value = '  +1.44%  '
$ws.Range("D3").Value = '1.805.31'
$ws.Range("E3").Value = '  +1.42%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.65'
$ws.Range("E5").Value = '  +0.71%  '
$ws.Range("E6").Value = '  +1.78%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '32.94'
$ws.Range("E8").Value = '  +3.87%  '
$ws.Range("E9").Value = '  +1.87%  '
$ws.Range("E10").Value = '  +0.82%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0949'
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("D12").Value = '2.064.47'
$ws.Range("E12").Value = '  +1.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.25'
$ws.Range("E13").Value = '  +2.79%  '
$ws.Range("D14").Value = '1.799.09'
$ws.Range("E14").Value = '  +1.17%  '
$ws.Range("E15").Value = '  +2.74%  '
$ws.Range("D16").Value = '34.658.89'
$ws.Range("E16").Value = '  +1.52%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.34'
$ws.Range("E17").Value = '  +3.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.04'
$ws.Range("E18").Value = '  +1.65%  '
$ws.Range("D19").Value = '0.0₃0805'
$ws.Range("E19").Value = '  +0.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '247.83'
$ws.Range("E20").Value = '  +0.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.36'
$ws.Range("E21").Value = '  +3.35%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.998'
$ws.Range("E22").Value = '  -0.38%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.20'
$ws.Range("E23").Value = '  +2.70%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '169.83'
$ws.Range("E24").Value = '  +4.59%  '
$ws.Range("E25").Value = '  +1.46%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.36'
$ws.Range("E26").Value = '  +2.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.67'
$ws.Range("E27").Value = '  +2.33%  '
$ws.Range("E28").Value = '  +2.18%  '
$ws.Range("E29").Value = '  -0.25%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.14'
$ws.Range("E30").Value = '  +11.63%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.25'
$ws.Range("E31").Value = '  +1.28%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0527'
$ws.Range("E32").Value = '  +1.20%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.82'
$ws.Range("E33").Value = '  +2.33%  '
$ws.Range("E34").Value = '  +2.94%  '
$ws.Range("D35").Value = '1.432.94'
$ws.Range("E35").Value = '  -0.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.61'
$ws.Range("E36").Value = '  +8.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.677'
$ws.Range("E37").Value = '  +2.89%  '
$ws.Range("E38").Value = '  +3.19%  '
$ws.Range("E39").Value = '  +0.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '85.39'
$ws.Range("E40").Value = '  +6.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.954'
$ws.Range("E41").Value = '  +3.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.39'
$ws.Range("E42").Value = '  +1.64%  '
$ws.Range("E43").Value = '  +3.50%  '
$ws.Range("E44").Value = '  +3.01%  '
$ws.Range("E45").Value = '  +2.77%  '
$ws.Range("E46").Value = '  +0.79%  '
$ws.Range("E47").Value = '  +0.60%  '
$ws.Range("D48").Value = '1.963.33'
$ws.Range("E48").Value = '  +1.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '105.98'
$ws.Range("E49").Value = '  +1.50%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.998'
$ws.Range("E50").Value = '  -0.32%  '
$ws.Range("E51").Value = '  -3.90%  '
